$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.8651844198272334
$ws1.Range("C2").Value = -0.8657304002748697
$ws1.Range("B3").Value = 1.334290919236224
$ws1.Range("C3").Value = 0.1253670301364903
$ws1.Range("B4").Value = 0.4436528822088864
$ws1.Range("C4").Value = -0.5174651542177248

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.446163129861008
$ws2.Range("C2").Value = -0.1681707613299935
$ws2.Range("B3").Value = 1.357860989289148
$ws2.Range("C3").Value = 0.3851991568590112
$ws2.Range("B4").Value = -0.7401770957417654
$ws2.Range("C4").Value = -0.1272403585205737
